$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (RE07 / UC07 "Enviar Boletos") previously listed "Sistema" as the actor
# and carried a long reviewer comment in column F. Update per commit:
# "Atualização - etapa 2 / inclusão do DSS UC05" -> actor corrected to
# "Gerente, Funcionário" and the now-resolved comment removed.
$ws.Range("E9").Value = "Gerente, Funcionário"
$ws.Range("F9").Value = ""

# With the long comment gone, the row's wrapped-text height shrinks from 75 to 45.
$ws.Rows(9).RowHeight = 45

# Update the active selection to F5, matching the saved view state.
$ws.Range("F5").Select()
